$wb = $excel.ActiveWorkbook
$busWs = $wb.Worksheets.Item("bus")
$branchWs = $wb.Worksheets.Item("branch")

# --- "bus" sheet: introduce a new "QL (pu)" column header/values ---
# The header row gains a new column; the two existing trailing headers
# shift one column to the right.
$busWs.Range("J1").Value = $busWs.Range("I1").Value()
$busWs.Range("I1").Value = $busWs.Range("H1").Value()
$busWs.Range("H1").Value = "QL (pu)"

# Row 2 gets a brand-new QL value (no shift needed - columns I/J were blank).
$busWs.Range("H2").Value = 0

# Row 4's old values move right to make room for the new QL value.
$busWs.Range("J4").Value = $busWs.Range("I4").Value()
$busWs.Range("I4").Value = $busWs.Range("H4").Value()
$busWs.Range("H4").Value = 0.4

# --- "branch" sheet: append two new branch rows ---
$branchWs.Range("A5").Value = 1
$branchWs.Range("B5").Value = 5
$branchWs.Range("C5").Value = 0.0015
$branchWs.Range("D5").Value = 0.02
$branchWs.Range("E5").Value = 0

$branchWs.Range("A6").Value = 3
$branchWs.Range("B6").Value = 4
$branchWs.Range("C6").Value = 0.00075
$branchWs.Range("D6").Value = 0.01
$branchWs.Range("E6").Value = 0

# --- selection / active sheet bookkeeping ---
[void]$busWs.Range("B4").Select()
$branchWs.Activate()
[void]$branchWs.Range("E6").Select()
